$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 998.625
$ws.Range("I12").Value = 298
$ws.Range("K12").Value = 298
$ws.Range("M12").Value = -128

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2053950
$ws.Range("J43").Value = 4100000
$ws.Range("L43").Value = 4100000
$ws.Range("N43").Value = -4100138

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 67268.17999999999
$ws.Range("J51").Value = 94427.71000000001
$ws.Range("L51").Value = 94427.71000000001
$ws.Range("N51").Value = -95395.71000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 3702.2222
$ws.Range("I53").Value = 2732.75
$ws.Range("K53").Value = 2732.75
$ws.Range("M53").Value = -2095.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 11409345
$ws.Range("I80").Value = 13900311
$ws.Range("J80").Value = 200000
$ws.Range("K80").Value = 41700933
$ws.Range("L80").Value = 600000
$ws.Range("M80").Value = -41699935
$ws.Range("N80").Value = -601996

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 11409345
$ws.Range("I83").Value = 13900311
$ws.Range("J83").Value = 200000
$ws.Range("K83").Value = 125102799
$ws.Range("L83").Value = 1800000
$ws.Range("M83").Value = -125097807
$ws.Range("N83").Value = -1809984

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 27838032
$ws.Range("J88").Value = 79271.22
$ws.Range("L88").Value = 79271.22
$ws.Range("N88").Value = -80083.22

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 27838032
$ws.Range("J91").Value = 79271.22
$ws.Range("L91").Value = 79271.22
$ws.Range("N91").Value = -82079.22

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3862.5454
$ws.Range("I137").Value = 5949.75
$ws.Range("K137").Value = 17849.25
$ws.Range("M137").Value = -15299.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5247.59
$ws.Range("I138").Value = 1528.3
$ws.Range("K138").Value = 4584.9
$ws.Range("M138").Value = 555.1000000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1986.5938
$ws.Range("I2").Value = 1145.5294
$ws.Range("K2").Value = 1145.5294
$ws.Range("M2").Value = -1032.5294

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2317773.8
$ws.Range("I32").Value = 2453716
$ws.Range("K32").Value = 2453716
$ws.Range("M32").Value = -2453429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4857.604
$ws.Range("I61").Value = 2316.0527
$ws.Range("J61").Value = 11296.2
$ws.Range("K61").Value = 2316.0527
$ws.Range("L61").Value = 11296.2
$ws.Range("M61").Value = -2104.0527
$ws.Range("N61").Value = -11720.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 41337.57
$ws.Range("I74").Value = 60068.355
$ws.Range("J74").Value = 3876
$ws.Range("K74").Value = 60068.355
$ws.Range("L74").Value = 3876
$ws.Range("M74").Value = -59194.355
$ws.Range("N74").Value = -5624

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 41337.57
$ws.Range("I77").Value = 60068.355
$ws.Range("J77").Value = 3876
$ws.Range("K77").Value = 300341.775
$ws.Range("L77").Value = 19380
$ws.Range("M77").Value = -295973.775
$ws.Range("N77").Value = -28116

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 8345348
$ws.Range("J97").Value = 16689469
$ws.Range("L97").Value = 16689469
$ws.Range("N97").Value = -16690461

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1986.5938
$ws.Range("I116").Value = 1145.5294
$ws.Range("K116").Value = 1145.5294
$ws.Range("M116").Value = 1148.4706

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4760.2656
$ws.Range("I132").Value = 3546.1592
$ws.Range("J132").Value = 7431.3
$ws.Range("K132").Value = 10638.4776
$ws.Range("L132").Value = 22293.9
$ws.Range("M132").Value = -8108.4776
$ws.Range("N132").Value = -27353.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4857.604
$ws.Range("I136").Value = 2316.0527
$ws.Range("J136").Value = 11296.2
$ws.Range("K136").Value = 6948.158100000001
$ws.Range("L136").Value = 33888.60000000001
$ws.Range("M136").Value = -4398.158100000001
$ws.Range("N136").Value = -38988.60000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1986.5938
$ws.Range("I3").Value = 1145.5294
$ws.Range("K3").Value = 1145.5294
$ws.Range("M3").Value = -1031.5294

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9581.828
$ws.Range("I31").Value = 3595.8333
$ws.Range("J31").Value = 12704.956
$ws.Range("K31").Value = 3595.8333
$ws.Range("L31").Value = 12704.956
$ws.Range("M31").Value = -3300.8333
$ws.Range("N31").Value = -13294.956

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 9581.828
$ws.Range("I34").Value = 3595.8333
$ws.Range("J34").Value = 12704.956
$ws.Range("K34").Value = 3595.8333
$ws.Range("L34").Value = 12704.956
$ws.Range("M34").Value = -3393.8333
$ws.Range("N34").Value = -13108.956

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 11454.454
$ws.Range("I99").Value = 17499.75
$ws.Range("J99").Value = 8000
$ws.Range("K99").Value = 17499.75
$ws.Range("L99").Value = 8000
$ws.Range("M99").Value = -16001.75
$ws.Range("N99").Value = -10996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 11454.454
$ws.Range("I126").Value = 17499.75
$ws.Range("J126").Value = 8000
$ws.Range("K126").Value = 52499.25
$ws.Range("L126").Value = 24000
$ws.Range("M126").Value = -50029.25
$ws.Range("N126").Value = -28940

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 5504.4
$ws.Range("I68").Value = 1299.75
$ws.Range("J68").Value = 22323
$ws.Range("K68").Value = 3899.25
$ws.Range("L68").Value = 66969
$ws.Range("M68").Value = -3088.25
$ws.Range("N68").Value = -68591

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 3242
$ws.Range("I69").Value = 3000
$ws.Range("J69").Value = 3968
$ws.Range("K69").Value = 9000
$ws.Range("L69").Value = 11904
$ws.Range("M69").Value = -8189
$ws.Range("N69").Value = -13526

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 5504.4
$ws.Range("I71").Value = 1299.75
$ws.Range("J71").Value = 22323
$ws.Range("K71").Value = 11697.75
$ws.Range("L71").Value = 200907
$ws.Range("M71").Value = -7641.75
$ws.Range("N71").Value = -209019

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H72").Value = 3242
$ws.Range("I72").Value = 3000
$ws.Range("J72").Value = 3968
$ws.Range("K72").Value = 27000
$ws.Range("L72").Value = 35712
$ws.Range("M72").Value = -22944
$ws.Range("N72").Value = -43824

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 166669090
$ws.Range("I80").Value = 166668670
$ws.Range("J80").Value = 166669500
$ws.Range("K80").Value = 500006010
$ws.Range("L80").Value = 500008500
$ws.Range("M80").Value = -500005074
$ws.Range("N80").Value = -500010372

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 166669090
$ws.Range("I83").Value = 166668670
$ws.Range("J83").Value = 166669500
$ws.Range("K83").Value = 1500018030
$ws.Range("L83").Value = 1500025500
$ws.Range("M83").Value = -1500013350
$ws.Range("N83").Value = -1500034860

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1062.7858
$ws.Range("I98").Value = 1059.8572
$ws.Range("J98").Value = 1065.7142
$ws.Range("K98").Value = 3179.5716
$ws.Range("L98").Value = 3197.1426
$ws.Range("M98").Value = -1681.5716
$ws.Range("N98").Value = -6193.142599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4475.625
$ws.Range("I80").Value = 3980.8
$ws.Range("K80").Value = 3980.8
$ws.Range("M80").Value = -2982.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4475.625
$ws.Range("I83").Value = 3980.8
$ws.Range("K83").Value = 19904
$ws.Range("M83").Value = -14912

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6998.577
$ws.Range("I7").Value = 6139.7144
$ws.Range("J7").Value = 8000.5835
$ws.Range("K7").Value = 6139.7144
$ws.Range("L7").Value = 8000.5835
$ws.Range("M7").Value = -6027.7144
$ws.Range("N7").Value = -8224.583500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4653
$ws.Range("I68").Value = 3633.7144
$ws.Range("K68").Value = 3633.7144
$ws.Range("M68").Value = -2884.7144

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 4653
$ws.Range("I71").Value = 3633.7144
$ws.Range("K71").Value = 18168.572
$ws.Range("M71").Value = -14424.572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3634.24
$ws.Range("I100").Value = 3131.8
$ws.Range("J100").Value = 3969.2
$ws.Range("K100").Value = 3131.8
$ws.Range("L100").Value = 3969.2
$ws.Range("M100").Value = -2590.8
$ws.Range("N100").Value = -5051.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 6998.577
$ws.Range("I126").Value = 6139.7144
$ws.Range("J126").Value = 8000.5835
$ws.Range("K126").Value = 18419.1432
$ws.Range("L126").Value = 24001.7505
$ws.Range("M126").Value = -15949.1432
$ws.Range("N126").Value = -28941.7505

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H135").Value = 85000
$ws.Range("J135").Value = 85000
$ws.Range("L135").Value = 85000
$ws.Range("N135").Value = -95140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 111587
$ws.Range("I122").Value = 244395.77
$ws.Range("J122").Value = 4075.1428
$ws.Range("K122").Value = 733187.3099999999
$ws.Range("L122").Value = 12225.4284
$ws.Range("M122").Value = -730737.3099999999
$ws.Range("N122").Value = -17125.4284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 6499.6665
$ws.Range("I126").Value = 1999
$ws.Range("J126").Value = 7399.8
$ws.Range("K126").Value = 5997
$ws.Range("L126").Value = 22199.4
$ws.Range("M126").Value = -3527
$ws.Range("N126").Value = -27139.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 14298102
$ws.Range("I132").Value = 20005206
$ws.Range("K132").Value = 60015618
$ws.Range("M132").Value = -60013088
